$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Putz 203")
Write-Output $ws.Tab.Color
$ws.Tab.Color = 5296274
Write-Output $ws.Tab.Color
